# Terraform DevOpsAgent.xlsx - "Working DevOps Agent Deployment"
#
# Replace the placeholder / stale variable values on the "Variables" sheet
# with the new working values, and update the active selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# aws_access_key (row 2)
$ws.Range("B2").Value = "xxxxxxxxxxxxxxxxxxxx"

# aws_secret_key (row 3)
$ws.Range("B3").Value = "xxxxxxxxxxxxxxxxxxxxxxxxxxxxxxxxxxxxxxxx"

# aws_region (row 4)
$ws.Range("B4").Value = "eu-west-2"

# corporate_ip (row 8)
$ws.Range("B8").Value = "212.139.37.214"

# instance_key (row 9) - was previously blank
$ws.Range("B9").Value = "<Key_Name>"

# devops_organisation (row 10)
$ws.Range("B10").Value = "<Org_Name>"

# devops_pat (row 11)
$ws.Range("B11").Value = "<Personal_Access_Token>"

# devops_pool_name (row 12)
$ws.Range("B12").Value = "<Pool_Name>"

# environment_tag (row 13)
$ws.Range("B13").Value = "Development"

# owner_tag (row 14)
$ws.Range("B14").Value = "<Owner>"

# Update the sheet's active selection to match the saved workbook state.
[void]$ws.Range("B18").Select()
